# Weekly update: insert two new "Choclo" price rows (newest report, date
# serial 44931) at the top of the Comercializadora del Agro de Limari /
# Choclo block (rows 137-138), pushing the existing rows 137-148 down to
# 139-150.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 137 (formatting/style of the
# row above - including the date format on column D - is carried down
# automatically by Insert).
$ws.Rows("137:138").Insert()

# New row 137: Choclo / Choclero / Primera
$ws.Range("A137").Value = 2
$ws.Range("B137").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C137").Value = "Coquimbo"
$ws.Range("D137").Value = 44931
$ws.Range("E137").Value = 4
$ws.Range("F137").Value = 100112024
$ws.Range("G137").Value = "Choclo"
$ws.Range("H137").Value = "Choclero"
$ws.Range("I137").Value = "Primera"
$ws.Range("J137").Value = 160000
$ws.Range("K137").Value = 250
$ws.Range("L137").Value = 300
$ws.Range("M137").Value = 275
$ws.Range("N137").Value = "$/unidad"
$ws.Range("O137").Value = "Provincia de Limarí"
$ws.Range("P137").Value = 275
$ws.Range("Q137").Value = 1
$ws.Range("R137").Value = "Hortaliza"

# New row 138: Choclo / Dulce o Americano / Primera
$ws.Range("A138").Value = 2
$ws.Range("B138").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C138").Value = "Coquimbo"
$ws.Range("D138").Value = 44931
$ws.Range("E138").Value = 4
$ws.Range("F138").Value = 100112024
$ws.Range("G138").Value = "Choclo"
$ws.Range("H138").Value = "Dulce o Americano"
$ws.Range("I138").Value = "Primera"
$ws.Range("J138").Value = 140000
$ws.Range("K138").Value = 90
$ws.Range("L138").Value = 100
$ws.Range("M138").Value = 95
$ws.Range("N138").Value = "$/unidad"
$ws.Range("O138").Value = "Provincia de Limarí"
$ws.Range("P138").Value = 95
$ws.Range("Q138").Value = 1
$ws.Range("R138").Value = "Hortaliza"
